$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns C (p_total), D (p_eligible), E (p_registered) for rows 2-53 ---
$cde = @(
    @(2, 5024279, 3716000, 2527000),
    @(3, 733391, 516000, 383000),
    @(4, 7151502, 5075000, 3878000),
    @(5, 3011524, 2195000, 1361000),
    @(6, 39538223, 25946000, 18001000),
    @(7, 5773714, 4200000, 4158895),
    @(8, 3605944, 2524000, 1850000),
    @(9, 989948, 722000, 542000),
    @(10, 689545, 534000, 464000),
    @(11, 21538187, 15645000, 14568993),
    @(12, 10711908, 7400000, 5233000),
    @(13, 1455271, 980000, 673000),
    @(14, 1839106, 1299000, 900000),
    @(15, 12812508, 8860000, 6590000),
    @(16, 6785528, 4921000, 3412000),
    @(17, 3190369, 2293000, 1742000),
    @(18, 2937880, 1975000, 1398000),
    @(19, 4505836, 3227000, 2450000),
    @(20, 4657757, 3299000, 2286000),
    @(21, 1362359, 1075000, 832000),
    @(22, 6177224, 4303000, 3383000),
    @(23, 7029917, 4897000, 4812909),
    @(24, 10077331, 7467000, 7151051),
    @(25, 5706494, 4142000, 3436000),
    @(26, 2961279, 2177000, 1749000),
    @(27, 6154913, 4475000, 3388000),
    @(28, 1084225, 827000, 641000),
    @(29, 1961504, 1369000, 971000),
    @(30, 3104614, 2198000, 1455000),
    @(31, 1377529, 1077000, 843000),
    @(32, 9288994, 5921000, 5008000),
    @(33, 2117522, 1498000, 1028000),
    @(34, 20201249, 13298000, 9370000),
    @(35, 10439388, 7391000, 7242242),
    @(36, 779094, 556000, 429000),
    @(37, 11799448, 8740000, 6733000),
    @(38, 3959353, 2800000, 1884000),
    @(39, 4237256, 3242000, 2590000),
    @(40, 13011844, 9621000, 7337000),
    @(41, 1097379, 776000, 575000),
    @(42, 5118425, 3878000, 2713000),
    @(43, 886667, 649000, 437000),
    @(44, 6910840, 5038000, 3742000),
    @(45, 29145505, 18581000, 13343000),
    @(46, 3271616, 2178000, 1681844),
    @(47, 643077, 500000, 495267),
    @(48, 8631393, 5974000, 4541000),
    @(49, 7705281, 5389000, 4892871),
    @(50, 1793716, 1379000, 928000),
    @(51, 5893718, 4421000, 3391000),
    @(52, 576851, 427000, 296000),
    @(53, 331458425, 231591000, 179737072)
)

foreach ($row in $cde) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "0"
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 5).NumberFormat = "0"
}

# --- Column B (ec) corrections ---
$bChanges = @(
    @(6, 54),
    @(7, 10),
    @(11, 30),
    @(15, 19),
    @(24, 15),
    @(28, 4),
    @(34, 28),
    @(35, 16),
    @(37, 17),
    @(39, 8),
    @(40, 19),
    @(45, 40),
    @(50, 4)
)
foreach ($row in $bChanges) {
    $ws.Cells.Item($row[0], 2).Value = $row[1]
}

# --- Column F (Republican) corrections ---
$fChanges = @(
    @(8, 714717),
    @(16, 1729519),
    @(31, 365660),
    @(34, 3244798),
    @(35, 2758775),
    @(40, 3377674),
    @(44, 1852475),
    @(53, 74216154)
)
foreach ($row in $fChanges) {
    $ws.Cells.Item($row[0], 6).Value = $row[1]
}

# --- Column G (Democrat) corrections ---
$gChanges = @(
    @(8, 1080831),
    @(16, 1242416),
    @(31, 424937),
    @(34, 5230985),
    @(40, 3458229),
    @(44, 1143711),
    @(53, 81268924)
)
foreach ($row in $gChanges) {
    $ws.Cells.Item($row[0], 7).Value = $row[1]
}

# --- Column H (other) corrections / de-formularization (rows 2-52 + total row 53) ---
$hChanges = @(
    @(6, 384202),
    @(7, 88021),
    @(9, 7475),
    @(12, 64473),
    @(16, 61186),
    @(18, 30574),
    @(22, 75593),
    @(23, 81998),
    @(24, 85410),
    @(25, 76029),
    @(26, 17597),
    @(27, 54212),
    @(28, 15286),
    @(29, 24954),
    @(31, 15608),
    @(32, 57744),
    @(34, 119043),
    @(40, 79380),
    @(41, 10349),
    @(42, 36685),
    @(44, 57665),
    @(45, 165583),
    @(46, 62867),
    @(47, 11904),
    @(48, 84526),
    @(49, 133368),
    @(50, 13365),
    @(51, 56991),
    @(52, 9715),
    @(53, 2898325)
)
foreach ($row in $hChanges) {
    $ws.Cells.Item($row[0], 8).Value = $row[1]
}

# --- Worksheet view changes: zoom to 100%, selection moves to J1:R1048576 ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("J1:R1048576").Select()

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1
